# Update cryptos list (rows 2-51) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.697.96'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '1.647.98'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''212.97'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").Value = '''0.533'
$ws.Range("E6").Value = '  +4.09%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''23.09'
$ws.Range("E8").Value = '  -2.50%  '
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("D12").Value = '1.882.45'
$ws.Range("E12").Value = '  -0.29%  '
$ws.Range("D13").Value = '1.635.86'
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("D14").Value = '''4.04'
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("E15").Value = '  -1.48%  '
$ws.Range("D16").Value = '''64.16'
$ws.Range("E16").Value = '  -2.26%  '
$ws.Range("D17").Value = '27.698.47'
$ws.Range("E17").Value = '  +1.18%  '
$ws.Range("D18").Value = '''230.13'
$ws.Range("D19").Value = '''7.66'
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("D20").Value = '0.0₃0723'
$ws.Range("E20").Value = '  -0.60%  '
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("E22").Value = '  -1.37%  '
$ws.Range("D23").Value = '''10.09'
$ws.Range("E23").Value = '  +7.99%  '
$ws.Range("E24").Value = '  -2.62%  '
$ws.Range("D25").Value = '''149.15'
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("D26").Value = '''6.97'
$ws.Range("E27").Value = '  +1.22%  '
$ws.Range("E30").Value = '  -0.56%  '
$ws.Range("D31").Value = '''0.0485'
$ws.Range("E31").Value = '  -2.54%  '
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("E33").Value = '  +1.87%  '
$ws.Range("D34").Value = '1.435.08'
$ws.Range("E34").Value = '  -1.42%  '
$ws.Range("E35").Value = '  +2.44%  '
$ws.Range("E36").Value = '  -1.85%  '
$ws.Range("D37").Value = '''0.573'
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("E38").Value = '  -2.76%  '
$ws.Range("E39").Value = '  -0.90%  '
$ws.Range("D40").Value = '''0.898'
$ws.Range("E40").Value = '  +14.07%  '
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("D42").Value = '''0.999'
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").Value = '''5.59'
$ws.Range("E43").Value = '  +2.24%  '
$ws.Range("E44").Value = '  +2.29%  '
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("D46").Value = '''65.41'
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("D47").Value = '1.791.90'
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("E48").Value = '  -1.70%  '
$ws.Range("D49").Value = '''86.54'
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("D50").Value = '0.0₇0999'
$ws.Range("E50").Value = '  -6.16%  '
$ws.Range("E51").Value = '  -2.22%  '

# Row 28/29 swap: BinanceUSD and EthereumClassic traded positions
$ws.Range("B28").Value = 'BinanceUSD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '''15.63'
$ws.Range("E29").Value = '  -1.66%  '
